$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.532132
$ws.Range("N2").Value = 79.596396
$ws.Range("O2").Value = 0.3960736634233649
$ws.Range("P2").Value = 0.3960736634233648
$ws.Range("Q2").Value = 1.890874295288
$ws.Range("R2").Value = 17.017868657592
$ws.Range("S2").Value = 0.3960736634233649
$ws.Range("T2").Value = 0.3960736634233648

# Row 3 (Target cluster: FAPs)
$ws.Range("O3").Value = 0.2505213219764053
$ws.Range("P3").Value = 0.2505213219764053
$ws.Range("S3").Value = 0.2505213219764053
$ws.Range("T3").Value = 0.2505213219764053

# Row 4 (Target cluster: MuSCs)
$ws.Range("M4").Value = 23.67385
$ws.Range("N4").Value = 71.02154999999999
$ws.Range("O4").Value = 0.3534050146002298
$ws.Range("P4").Value = 0.3534050146002298
$ws.Range("Q4").Value = 1.687172159233333
$ws.Range("R4").Value = 15.1845494331
$ws.Range("S4").Value = 0.3534050146002298
$ws.Range("T4").Value = 0.3534050146002298
